$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "65.560.87"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +1.89%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.648.05"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +0.86%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "605.18"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +1.56%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "156.34"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +2.61%  "
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -0.10%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.588"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -0.31%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.646.34"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +0.87%  "
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +7.69%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.402"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +2.19%  "
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +0.61%  "
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +1.49%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "29.72"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +5.85%  "
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +13.97%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.126.25"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +0.95%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "65.299.81"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.648.21"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -0.94%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.66"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +2.84%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.88"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +2.30%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "358.72"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +2.57%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "7.45"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +5.22%  "
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -0.08%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "69.76"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +3.02%  "
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -0.14%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.42"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +2.24%  "
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +15.39%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.63"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -2.45%  "
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +1.99%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "8.08"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -3.20%  "
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -0.04%  "
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +4.78%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "524.46"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -5.31%  "
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -3.29%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.51"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -0.05%  "
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +1.73%  "
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +2.11%  "
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +2.82%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "161.87"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -2.51%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.97"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -0.94%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.999"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -0.01%  "
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +0.01%  "
$ws.Range("B43").NumberFormat = "@"
$ws.Range("B43").Value = "Aave"
$ws.Range("C43").NumberFormat = "@"
$ws.Range("C43").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "165.53"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -2.02%  "
$ws.Range("B44").NumberFormat = "@"
$ws.Range("B44").Value = "OKB"
$ws.Range("C44").NumberFormat = "@"
$ws.Range("C44").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "41.77"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +1.36%  "
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +0.65%  "
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +5.51%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0608"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +2.92%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "22.93"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -1.14%  "
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +1.67%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0262"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +3.63%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0979"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +0.26%  "
